$wb = $excel.ActiveWorkbook

$wsAlpha = $wb.Worksheets.Item("Alpha")
$wsAlpha.Range("H2").Value = -1.98908430042096
$wsAlpha.Range("H3").Value = -1.448638050480701

$wsDelta = $wb.Worksheets.Item("Delta")
$wsDelta.Range("H2").Value = -2.833048256440481
$wsDelta.Range("H3").Value = -2.395648834918044
